$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.670.29"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.585.80"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.40"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -3.11%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.27"
$ws.Range("E8").Value = "  -4.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.253"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0868"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "1.809.69"
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").Value = "1.584.11"
$ws.Range("E13").Value = "  -3.20%  "
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("E15").Value = "  -5.81%  "
$ws.Range("D16").Value = "27.645.68"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "219.36"
$ws.Range("E18").Value = "  -4.20%  "
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  -5.34%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.53"
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("E24").Value = "  -5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.89"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.10"
$ws.Range("E28").Value = "  -2.92%  "
$ws.Range("E29").Value = "  -3.91%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("D33").Value = "1.382.60"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.961"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.70"
$ws.Range("E44").Value = "  -3.90%  "
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  -4.12%  "
$ws.Range("D47").Value = "1.721.51"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.28"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("E50").Value = "  -5.11%  "
$ws.Range("E51").Value = "  -1.01%  "
